$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append after the last existing row (161)
$row162 = @(11, "Vega Monumental Concepción", "Bíobío", 44595, 8, 100112040, "Cilantro", "Sin especificar", "Primera", 200, 600, 700, 650, "$/atado 0,5 a 1 kilo", "Región de Ñuble", 650, 1, "Hortaliza")
$row163 = @(11, "Vega Monumental Concepción", "Bíobío", 44595, 8, 100112040, "Cilantro", "Sin especificar", "Segunda", 100, 500, 500, 500, "$/atado 0,5 a 1 kilo", "Región de Ñuble", 500, 1, "Hortaliza")

for ($c = 0; $c -lt $row162.Length; $c++) {
    $ws.Cells.Item(162, $c + 1).Value = $row162[$c]
}

for ($c = 0; $c -lt $row163.Length; $c++) {
    $ws.Cells.Item(163, $c + 1).Value = $row163[$c]
}

# Column D (4) holds a date value; match the existing date number format used in the sheet
$ws.Cells.Item(162, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(163, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
